$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell, forcing text storage even when
# the string looks like a plain number (prevents Excel from coercing it
# to a floating-point number), while leaving the cell style unchanged.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "28.739.69"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "1.575.04"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "213.60"
$ws.Range("E5").Value = "  +0.14%  "
Set-TextValue "D6" "0.492"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  +0.09%  "
Set-TextValue "D8" "44.71"
$ws.Range("E8").Value = "  +1.84%  "
Set-TextValue "D9" "24.16"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("E11").Value = "  -0.59%  "
Set-TextValue "D12" "0.0892"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "1.799.08"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "1.566.95"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "28.720.15"
$ws.Range("E16").Value = "  +1.73%  "
Set-TextValue "D18" "62.47"
$ws.Range("E18").Value = "  -1.04%  "
Set-TextValue "D19" "231.41"
$ws.Range("E19").Value = "  +1.97%  "
Set-TextValue "D20" "7.38"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").Value = "0.0₃0694"
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -4.44%  "
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  +6.49%  "
Set-TextValue "D26" "151.90"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  -0.49%  "
Set-TextValue "D34" "3.12"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "1.396.92"
$ws.Range("E35").Value = "  -0.21%  "
Set-TextValue "D36" "1.05"
$ws.Range("E36").Value = "  +2.38%  "
Set-TextValue "D37" "1.55"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -1.84%  "
Set-TextValue "D44" "1.90"
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("E46").Value = "  -1.76%  "
Set-TextValue "D47" "0.962"
$ws.Range("E47").Value = "  -1.93%  "
Set-TextValue "D48" "63.25"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").Value = "1.711.43"
$ws.Range("E49").Value = "  -0.69%  "
Set-TextValue "D50" "86.66"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  -0.45%  "
